# Auto-generated edit script applying numeric corrections to the Siren_Profits workbook
# (per-sheet leve profit calculations refreshed by the scheduled price-update runner)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5333.1665
$ws.Range("J40").Value = 4999.8
$ws.Range("L40").Value = 4999.8
$ws.Range("N40").Value = -5349.8
$ws.Range("H107").Value = 30000
$ws.Range("H129").Value = 100001700
$ws.Range("I129").Value = 1481.8334
$ws.Range("J129").Value = 250002030
$ws.Range("K129").Value = 4445.5002
$ws.Range("L129").Value = 750006090
$ws.Range("M129").Value = 554.4997999999996
$ws.Range("N129").Value = -750016090
$ws.Range("H137").Value = 10859.917
$ws.Range("I137").Value = 12248.3
$ws.Range("J137").Value = 3918
$ws.Range("K137").Value = 36744.89999999999
$ws.Range("L137").Value = 11754
$ws.Range("M137").Value = -34194.89999999999
$ws.Range("N137").Value = -16854
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3699.12
$ws.Range("I32").Value = 3699.12
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 3699.12
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -3412.12
$ws.Range("N32").ClearContents()
$ws.Range("H110").Value = 2879.25
$ws.Range("I110").Value = 1607.8
$ws.Range("K110").Value = 1607.8
$ws.Range("M110").Value = 437.2
$ws.Range("H132").Value = 2099.147
$ws.Range("I132").Value = 1075.4814
$ws.Range("K132").Value = 3226.4442
$ws.Range("M132").Value = -696.4441999999999
$ws.Range("H134").Value = 364285.16
$ws.Range("I134").Value = 299999
$ws.Range("J134").Value = 374999.5
$ws.Range("K134").Value = 299999
$ws.Range("L134").Value = 374999.5
$ws.Range("M134").Value = -294929
$ws.Range("N134").Value = -385139.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H16").Value = 8008135
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("H134").Value = 4312.7393
$ws.Range("I134").Value = 3031
$ws.Range("J134").Value = 7242.4287
$ws.Range("K134").Value = 9093
$ws.Range("L134").Value = 21727.2861
$ws.Range("M134").Value = -6558
$ws.Range("N134").Value = -26797.2861
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1598.25
$ws.Range("J16").Value = 1400
$ws.Range("L16").Value = 1400
$ws.Range("N16").Value = -1974
$ws.Range("H22").Value = 663.2353000000001
$ws.Range("I22").Value = 282
$ws.Range("J22").Value = 744.9286
$ws.Range("K22").Value = 282
$ws.Range("L22").Value = 744.9286
$ws.Range("M22").Value = 68
$ws.Range("N22").Value = -1444.9286
$ws.Range("H31").Value = 3903.3635
$ws.Range("I31").Value = 2789.75
$ws.Range("J31").Value = 5239.7
$ws.Range("K31").Value = 2789.75
$ws.Range("L31").Value = 5239.7
$ws.Range("M31").Value = -2494.75
$ws.Range("N31").Value = -5829.7
$ws.Range("H34").Value = 3903.3635
$ws.Range("I34").Value = 2789.75
$ws.Range("J34").Value = 5239.7
$ws.Range("K34").Value = 2789.75
$ws.Range("L34").Value = 5239.7
$ws.Range("M34").Value = -2587.75
$ws.Range("N34").Value = -5643.7
$ws.Range("H86").Value = 12624.25
$ws.Range("I86").Value = 9936.375
$ws.Range("K86").Value = 9936.375
$ws.Range("M86").Value = -8813.375
$ws.Range("H89").Value = 12624.25
$ws.Range("I89").Value = 9936.375
$ws.Range("K89").Value = 49681.875
$ws.Range("M89").Value = -44065.875
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()
$ws.Range("H109").Value = 75000
$ws.Range("J109").Value = 75000
$ws.Range("L109").Value = 75000
$ws.Range("N109").Value = -77080
$ws.Range("H113").Value = 1598.25
$ws.Range("J113").Value = 1400
$ws.Range("L113").Value = 1400
$ws.Range("N113").Value = -5740
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 2499.0715
$ws.Range("I51").Value = 998
$ws.Range("K51").Value = 2994
$ws.Range("M51").Value = -2534
$ws.Range("H56").Value = 8720.429
$ws.Range("I56").Value = 8720.429
$ws.Range("K56").Value = 8720.429
$ws.Range("M56").Value = -8190.429
$ws.Range("H86").Value = 1299.8572
$ws.Range("I86").Value = 1049.75
$ws.Range("K86").Value = 3149.25
$ws.Range("M86").Value = -1963.25
$ws.Range("H89").Value = 1299.8572
$ws.Range("I89").Value = 1049.75
$ws.Range("K89").Value = 9447.75
$ws.Range("M89").Value = -3519.75
$ws.Range("H92").Value = 250.66667
$ws.Range("I92").Value = 245.16667
$ws.Range("J92").Value = 261.66666
$ws.Range("K92").Value = 735.50001
$ws.Range("L92").Value = 784.9999799999999
$ws.Range("M92").Value = 512.49999
$ws.Range("N92").Value = -3280.99998
$ws.Range("H107").Value = 574.549
$ws.Range("I107").Value = 198.53334
$ws.Range("J107").Value = 731.2222
$ws.Range("K107").Value = 595.6000200000001
$ws.Range("L107").Value = 2193.6666
$ws.Range("M107").Value = 1324.39998
$ws.Range("N107").Value = -6033.6666
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 30623.666
$ws.Range("J122").Value = 23935.75
$ws.Range("L122").Value = 71807.25
$ws.Range("N122").Value = -76707.25
$ws.Range("H132").Value = 4173.8335
$ws.Range("I132").Value = 4122.364
$ws.Range("J132").Value = 4254.7144
$ws.Range("K132").Value = 12367.092
$ws.Range("L132").Value = 12764.1432
$ws.Range("M132").Value = -9837.091999999999
$ws.Range("N132").Value = -17824.1432
$ws.Range("H136").Value = 157581
$ws.Range("J136").Value = 157581
$ws.Range("L136").Value = 472743
$ws.Range("N136").Value = -477843
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3297
$ws.Range("I122").Value = 3013.25
$ws.Range("K122").Value = 9039.75
$ws.Range("M122").Value = -6589.75
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 78036
$ws.Range("I107").Value = 4048.1667
$ws.Range("K107").Value = 12144.5001
$ws.Range("M107").Value = -10224.5001
$ws.Range("H122").Value = 17235.158
$ws.Range("I122").Value = 3180.5862
$ws.Range("K122").Value = 9541.758600000001
$ws.Range("M122").Value = -7091.758600000001
$ws.Range("H132").Value = 9016.111000000001
$ws.Range("I132").Value = 9862.380999999999
$ws.Range("J132").Value = 6054.1665
$ws.Range("K132").Value = 29587.143
$ws.Range("L132").Value = 18162.4995
$ws.Range("M132").Value = -27057.143
$ws.Range("N132").Value = -23222.4995
$ws.Range("H138").Value = 80000
$ws.Range("J138").Value = 80000
$ws.Range("L138").Value = 80000
$ws.Range("N138").Value = -90280
